$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

$ws.Range("A9").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-18 for illustrative purposes only and are subject to change."

$ws.Range("D2").Value = 0.2492382959868233
$ws.Range("E2").Value = -0.005581631358619776

$ws.Range("D3").Value = 0.2482758731546902
$ws.Range("E3").Value = 0.005226480836237002

$ws.Range("D4").Value = 0.2506437559972614
$ws.Range("E4").Value = -0.02774869109947631

$ws.Range("D5").Value = 0.2518420748612253
$ws.Range("E5").Value = -0.01657458563535885

$ws.Range("E6").Value = -0.01122276139306022
